$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row relabeling
$ws.Range("E1").Value = "Весовая_категория"
$ws.Range("G1").Value = "Побед|Встреч"

# New column width for E (weight category column)
$ws.Columns("E").ColumnWidth = 19.5

# Update selection
$ws.Range("G1").Select()
